$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E values (Wn %/Physical property column) for rows 4-8
$ws.Range("E4").Value = 1.7
$ws.Range("E5").Value = 1.7
$ws.Range("E6").Value = 1.8
$ws.Range("E7").Value = 1.8
$ws.Range("E8").Value = 1.8

# Update the active cell selection to G4
$ws.Range("G4").Select()
